$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("E2").Value = 2
$ws.Range("A3").Value = 2
$ws.Range("E3").Value = 3
$ws.Range("A4").Value = 5
$ws.Range("E4").Value = 6
$ws.Range("A5").Value = 8
$ws.Range("E5").Value = 9
$ws.Range("A6").Value = 10
$ws.Range("E6").Value = 11
$ws.Range("A7").Value = 11
$ws.Range("E7").Value = 12
$ws.Range("A8").Value = 13
$ws.Range("E8").Value = 14
$ws.Range("A9").Value = 14
$ws.Range("E9").Value = 16
$ws.Range("A10").Value = 16
$ws.Range("E10").Value = 18
$ws.Range("A11").Value = 18
$ws.Range("E11").Value = 20
$ws.Range("A12").Value = 21
$ws.Range("E12").Value = 23
$ws.Range("A13").Value = 22
$ws.Range("E13").Value = 24
$ws.Range("A14").Value = 25
$ws.Range("E14").Value = 26
$ws.Range("A15").Value = 18
$ws.Range("E15").Value = 13
$ws.Range("A16").Value = 26
$ws.Range("E16").Value = 19
$ws.Range("A17").Value = 29
$ws.Range("E17").Value = 21
$ws.Range("A18").Value = 15
$ws.Range("E18").Value = 8
$ws.Range("A19").Value = 6
$ws.Range("E19").Value = 4
$ws.Range("A20").Value = 5
$ws.Range("E20").Value = 1
$ws.Range("A21").Value = 9
$ws.Range("E21").Value = 7
$ws.Range("A22").Value = 14
$ws.Range("E22").Value = 10
$ws.Range("A23").Value = 23
$ws.Range("E23").Value = 17
$ws.Range("A24").Value = 30
$ws.Range("E24").Value = 25
$ws.Range("A25").Value = 21
$ws.Range("E25").Value = 22
$ws.Range("A26").Value = 4
$ws.Range("E26").Value = 5
$ws.Range("A27").Value = 35
$ws.Range("E27").Value = 15
